# Regenerate the "K" column (G) of the save-data sheet with freshly
# calculated strikeout (K) values, replacing the old Strike# counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number => new K value (column G), rows 2..57
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 2
    10 = 1
    11 = 3
    12 = 4
    13 = 0
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 4
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 1
    31 = 0
    32 = 1
    33 = 2
    34 = 0
    35 = 1
    36 = 0
    37 = 2
    38 = 1
    39 = 2
    40 = 1
    41 = 1
    42 = 0
    43 = 0
    44 = 2
    45 = 1
    46 = 0
    47 = 1
    48 = 3
    49 = 1
    50 = 3
    51 = 2
    52 = 1
    53 = 1
    54 = 2
    55 = 1
    56 = 1
    57 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
